$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper pattern: set NumberFormat to Text before writing a
# numeric-looking string, then restore the Normal style so no visible
# format/style change is left behind (matches original plain inline-string cells).

$ws.Range("D2").Value = '69.989.20'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '3.695.48'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '648.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +1.54%  '

$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.19%  '

$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("D13").Value = '4.318.07'
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.01%  '

$ws.Range("D15").Value = '3.689.31'
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").Value = '69.920.73'
$ws.Range("E16").Value = '  +0.53%  '

$ws.Range("E17").Value = '  +0.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("E19").Value = '  +0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '472.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.654'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.34%  '

$ws.Range("D24").Value = '3.840.45'
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("E25").Value = '  +2.08%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  +1.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("E29").Value = '  -1.60%  '

$ws.Range("E30").Value = '  -1.69%  '

$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("E32").Value = '  -0.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.53%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.167'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.49%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.46%  '

$ws.Range("D36").Value = '3.691.39'
$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.64%  '

$ws.Range("E40").Value = '  +0.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '178.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0907'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.934'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("E51").Value = '  -3.29%  '
